$wb = $excel.ActiveWorkbook

# Sheet 2 = 建物 (Building): fix property_category from "land" to "building" for rows 2-7, column I
$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 7; $r++) {
    $wsBuilding.Range("I$r").Value = "building"
}

# Sheet 3 = 汽車 (Car): fix property_category from "land" to "car" for row 2, column H
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
